$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Site Code" column (A) for the last two data rows held text values
# "NB11" / "NB13" (site codes from a different agency numbering scheme).
# Re-importing the data replaced these with their corresponding plain
# numeric site codes (11 / 13), matching the numeric codes used by every
# other row in the column. Their "Agency Code" (D) stays "SW" - Excel
# automatically drops the now-unused "NB11"/"NB13" shared strings and
# renumbers the shared-string table accordingly.
$ws.Range("A16").Value = 11
$ws.Range("A17").Value = 13

# Reflect the final active selection left by the import/plotting session.
$ws.Range("A17").Select()
